$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (column D) values; these are numeric-looking strings
# that must remain stored as text, matching the original inlineStr cells. ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "243.50"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.05"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.394"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05937"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.402"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8097"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9113"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07429"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03322"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03066"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.940"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001586"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04812"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0005946"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.005556"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.004423"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0009872"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.00007804"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.632"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.434"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.150"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1349"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03870"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006216"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1065"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002901"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.006579"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005178"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005806"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.8351"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002264"

# --- Update Coin / Link / Volume(1h) text columns (rows 17-24 shifted) ---
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E17").Value = "16OneONE"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("E21").Value = "20NitroExNTXBestin24h"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("B23").Value = "KuCoinToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("E23").Value = "22KuCoinTokenKCS"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("E24").Value = "23BTSETokenBTSE"
